# Apply the commit's edits:
#  1. Refresh the cached "datetimeFigureOut" date field text (11.10.2024 -> 19.10.2024)
#     on the slide master and on every slide layout (the Date Placeholder shapes).
#  2. Rename the "Sales Analysis" slide title to "Sales Performance Analysis".
#  3. Rename the "Product Analysis" slide title to "Product Performance Analysis".

$p = $ppt.ActivePresentation

$oldDate = "11.10.2024"
$newDate = "19.10.2024"

# --- 1a. Slide master's Date Placeholder ---------------------------------
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "*Date*") {
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Every slide layout's Date Placeholder ----------------------------
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "*Date*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Slide 2 title: "Sales Analysis" -> "Sales Performance Analysis" ---
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item("Title")
if ($title2.TextFrame.TextRange.Text -eq "Sales Analysis") {
    $title2.TextFrame.TextRange.Text = "Sales Performance Analysis"
}

# --- 3. Slide 3 title: "Product Analysis" -> "Product Performance Analysis"
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item("Title")
if ($title3.TextFrame.TextRange.Text -eq "Product Analysis") {
    $title3.TextFrame.TextRange.Text = "Product Performance Analysis"
}
